# Task List For RPG Hero.xlsx - apply weekly update
#
# 1. "Tasks 02-11 to 02-18" sheet (currently 2nd tab) gets two extra
#    finished/in-progress tasks appended before its totals row, and two
#    "Over/Under" cells filled in.
# 2. A brand new "Tasks 02-18 to 02-25" sheet is inserted right after
#    "TaskList" with the next week's task list.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "Tasks 02-11 to 02-18" sheet (will shift to 3rd position
#    once the new sheet is inserted in front of it).
# ---------------------------------------------------------------------
$wsPrev = $wb.Worksheets.Item("Tasks 02-11 to 02-18")

# Fill in the two previously-blank "Over/Under" cells.
$wsPrev.Range("D12").Value = 0
$wsPrev.Range("D13").Value = 0

# Make room for two new task rows right before the blank separator /
# totals rows (old rows 14 & 15 -> new rows 16 & 17).
$wsPrev.Rows("14:15").Insert()

# New row 14: "Destroy Icons after leaving the camera" - Done.
$wsPrev.Range("A14").Value = "Destroy Icons after leaving the camera"
$wsPrev.Range("B14").Value = 0.25
$wsPrev.Range("C14").Value = 0.25
$wsPrev.Range("D14").Value = 0
$wsPrev.Range("E14").Value = "James"
$wsPrev.Range("F14").Value = "Done"
$wsPrev.Range("F13").Copy()
$wsPrev.Range("F14").PasteSpecial(-4122)

# New row 15: "Have Multiple methods called on certain situations" - In Progress.
$wsPrev.Range("A15").Value = "Have Multiple methods called on certain situations"
$wsPrev.Range("B15").Value = 4
$wsPrev.Range("C15").Value = 3
$wsPrev.Range("E15").Value = "James"
$wsPrev.Range("F15").Value = "In Progress"
$wsPrev.Range("G15").Value = "Issue with changing state "
# Row 15 inherited the "Done" fill from row 13 on insert - restyle to
# the "In Progress" look (copy format from the legend's In Progress cell).
$wsPrev.Range("F24").Copy()
$wsPrev.Range("F15").PasteSpecial(-4122)

# Totals row moved from 15 to 17; widen the summed ranges to include the
# two new task rows.
$wsPrev.Range("B17").Formula = "=SUM(B2:B15)"
$wsPrev.Range("C17").Formula = "=SUM(C2:C16)"

$wsPrev.Activate()
$wsPrev.Range("D16").Select()

# ---------------------------------------------------------------------
# 2. Insert the new "Tasks 02-18 to 02-25" sheet right after "TaskList".
# ---------------------------------------------------------------------
$wsTaskList = $wb.Worksheets.Item("TaskList")
$wsNew = $wb.Worksheets.Add($null, $wsTaskList)
$wsNew.Name = "Tasks 02-18 to 02-25"

# Header row.
$wsNew.Range("A1").Value = "Task"
$wsNew.Range("B1").Value = "Time Estimated to Complete"
$wsNew.Range("C1").Value = "Time Spent"
$wsNew.Range("D1").Value = "Over/Under"
$wsNew.Range("E1").Value = "Assigned To"
$wsNew.Range("F1").Value = "Status"
$wsNew.Range("G1").Value = "Notes"
$wsPrev.Range("A1:G1").Copy()
$wsNew.Range("A1:G1").PasteSpecial(-4122)

# Row 2: Expand MagicIcon class.
$wsNew.Range("A2").Value = "Expand MagicIcon class "
$wsNew.Range("B2").Value = 3
$wsNew.Range("E2").Value = "James"
$wsNew.Range("F2").Value = "TODO"
$wsNew.Range("G2").Value = "Implement Fire and Ice Spells"

# Row 3: Create StandingEnemy Class.
$wsNew.Range("A3").Value = "Create StandingEnemy Class"
$wsNew.Range("B3").Value = 2
$wsNew.Range("E3").Value = "James"
$wsNew.Range("F3").Value = "TODO"
$wsNew.Range("G3").Value = "Will be like a wolf "

# Row 4: Create FlyingEnemy Class.
$wsNew.Range("A4").Value = "Create FlyingEnemy Class"
$wsNew.Range("B4").Value = 4
$wsNew.Range("E4").Value = "James"
$wsNew.Range("F4").Value = "TODO"
$wsNew.Range("G4").Value = "Like a Eagle"

# Row 5: Create Base Boss Class.
$wsNew.Range("A5").Value = "Create Base Boss Class"
$wsNew.Range("B5").Value = 3
$wsNew.Range("E5").Value = "James"
$wsNew.Range("F5").Value = "TODO"
$wsNew.Range("G5").Value = "Will have a special attack"

# Row 6: Create ComparedItem Game Object.
$wsNew.Range("A6").Value = "Create ComparedItem Game Object"
$wsNew.Range("B6").Value = 2
$wsNew.Range("E6").Value = "James"
$wsNew.Range("F6").Value = "TODO"
$wsNew.Range("G6").Value = "Still need to alter to show compared magic when looking at magic section"

# Row 7: Research and implement Finite state machine for enemy.
$wsNew.Range("A7").Value = "Research and implement Finite state machine for enemy"
$wsNew.Range("B7").Value = 1
$wsNew.Range("E7").Value = "James"
$wsNew.Range("F7").Value = "In Progress"
$wsNew.Range("G7").Value = "Fix any troubles discovered during meeting"

# Row 8: Have Multiple methods called on certain situations.
$wsNew.Range("A8").Value = "Have Multiple methods called on certain situations"
$wsNew.Range("B8").Value = 2
$wsNew.Range("E8").Value = "James"
$wsNew.Range("F8").Value = "In Progress"
$wsNew.Range("G8").Value = "Issue with changing state "

# Status-cell colouring: copy the TODO / In Progress look from the
# "Tasks 02-11 to 02-18" sheet's legend cells.
$wsPrev.Range("F23").Copy()
$wsNew.Range("F2:F6").PasteSpecial(-4122)
$wsPrev.Range("F24").Copy()
$wsNew.Range("F7:F8").PasteSpecial(-4122)

# Totals row.
$wsNew.Range("A10").Value = "Total Hours Assigned"
$wsNew.Range("B10").Formula = "=SUM(B2:B8)"
$wsNew.Range("C10").Formula = "=SUM(C2:C9)"

# Legend block (TODO / In Progress / Done), copied from the previous
# sheet's own legend so the fills match exactly.
$wsPrev.Range("F23").Copy()
$wsNew.Range("L16").PasteSpecial(-4122)
$wsNew.Range("L16").Value = "TODO"
$wsPrev.Range("F24").Copy()
$wsNew.Range("L17").PasteSpecial(-4122)
$wsNew.Range("L17").Value = "In Progress"
$wsPrev.Range("F25").Copy()
$wsNew.Range("L18").PasteSpecial(-4122)
$wsNew.Range("L18").Value = "Done"

$wsNew.Activate()
$wsNew.Range("A8").Select()
